# Work on Scala, introduced type hierarchy and partial functions, and cleaned up the rest
#
# The backlog table (Table1 on Sheet1) gets a new task row inserted right
# after the "Bibliography" row (priority 5), and two existing Scala tasks
# get marked as DONE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (pushes the existing rows 9..21 down to 10..22)
$ws.Rows.Item(9).Insert()

# New backlog entry: Scala / rework the general structure / priority 5 / DONE
$ws.Range("A9").Value = "Scala"
$ws.Range("B9").Value = "rework the general structure"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = "DONE"

# Mark the existing Scala tasks as completed
$ws.Range("D12").Value = "DONE"   # Scala / type system with hierarchy
$ws.Range("D18").Value = "DONE"   # Scala / partial functions

# Grow the table to include the newly inserted row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D22"))

# Restore the sort (by Priority, then Section) across the new range
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("C2:C22"))
$tbl.Sort.SortFields.Add($ws.Range("A2:A22"))
$tbl.Sort.Apply()

# Leave the selection where the author left it
$ws.Range("D15").Select()
